$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores these cells as text (inline strings), even
# though several of the new "Price"/column-D values look numeric (e.g.
# "581.05"). Force NumberFormat to Text ("@") before assigning so Excel
# doesn't silently coerce them to numbers and lose the original formatting.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.698.77'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.271.12'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.05'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.78'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.66%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.36%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.409'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.835.91'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.40'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -5.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.732.98'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.269.43'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.44'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '401.53'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.56'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.95'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000118'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.187'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.57%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.65'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.47'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.92'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.54%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '163.92'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.45'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.23%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.87'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.804'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.41%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.35'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.679.23'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.92%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.45'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.47%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.72'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0678'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.59'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '334.77'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0273'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.28'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.966'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.50%  '
